$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 194, shifting existing rows 194..240 down to 195..241.
$ws.Rows.Item(194).Insert()

# Populate the new row 194 with the new weekly record (matches the
# constant columns used throughout this sheet, plus the new date/volume).
$ws.Range("A194").Value = 3
$ws.Range("B194").Value = "Femacal de La Calera"
$ws.Range("C194").Value = "Coquimbo"
$ws.Range("D194").Value = 44932
$ws.Range("E194").Value = 5
$ws.Range("F194").Value = 100112010
$ws.Range("G194").Value = "Achicoria"
$ws.Range("H194").Value = "Sin especificar"
$ws.Range("I194").Value = "Primera"
$ws.Range("J194").Value = 65
$ws.Range("K194").Value = 7000
$ws.Range("L194").Value = 7000
$ws.Range("M194").Value = 7000
$ws.Range("N194").Value = "$/caja 16 unidades"
$ws.Range("O194").Value = "Provincia de Quillota"
$ws.Range("P194").Value = 438
$ws.Range("Q194").Value = 16
$ws.Range("R194").Value = "Hortaliza"

# Match the date-format style used by the rest of column D.
$ws.Range("D194").NumberFormat = $ws.Range("D195").NumberFormat
